$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'316.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.69%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.73%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.122"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.53%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08168"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.89%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.984"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.92%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.368"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.34%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'2.21%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.99%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'-6.51%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1972"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.96%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'-0.58%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03505"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.32%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09755"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.67%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001407"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.18%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006058"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.77%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.649"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-7.52%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E19").Value = "'1.81%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1316"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.44%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.966"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'7.13%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2492"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.77%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04378"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.05%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001241"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.15%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004759"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'9.70%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003896"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'199.55%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-7.66%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02209"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'8.99%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05190"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.78%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007764"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.76%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01026"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.86%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1401"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.85%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002103"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.30%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.47%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006908"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'8.48%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'0.002885"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'0.65%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'30.13%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
